# Swap the two theme color palettes that live in ppt/theme/theme1.xml and
# ppt/theme/theme2.xml. Before the edit: theme1.xml held the stock
# "Office Theme" palette and theme2.xml held the "Integral" / "Red Violet"
# palette that the deck's slide master actually uses. After the edit the
# palettes trade places: the slide master's theme should show the stock
# Office colours, as if the two files' contents were exchanged.

$p = $ppt.ActivePresentation

# The 12 slots of a DrawingML colour scheme, in Office/ECMA order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# These are the values that originally lived in ppt/theme/theme1.xml
# ("Office Theme"); applying them to the live theme reproduces the swap.
$officeThemeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $tcs.Colors($i).RGB = $officeThemeColors[$i - 1]
}
